$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-14"

# Update the header label text (shared string) in I1
$ws.Range("I1").Value = "2022 (through 03-14)"

# Update March total (row 4) for column I
$ws.Range("I4").Value = 65

# Update the overall Total row (row 14) for column I
$ws.Range("I14").Value = 365
